$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: reduce stored precision to 2 decimal places ("custom accuracy") ---
$ws.Range("B5").Value = 9.67
$ws.Range("C5").Value = 6.96
$ws.Range("D5").Value = 0.19
$ws.Range("E5").Value = 19.19
$ws.Range("F5").Value = 16.23
$ws.Range("G5").Value = 7.46
$ws.Range("H5").Value = 26.68
$ws.Range("I5").Value = 10.67
$ws.Range("J5").Value = 4.93
$ws.Range("K5").Value = 7.81
$ws.Range("L5").Value = 7.82
$ws.Range("M5").Value = 8.04
$ws.Range("N5").Value = 2.3
$ws.Range("O5").Value = 6.9
$ws.Range("P5").Value = 10.45
$ws.Range("Q5").Value = 5.77
$ws.Range("R5").Value = 0.31
$ws.Range("S5").Value = 0.18
$ws.Range("T5").Value = 101.19
$ws.Range("U5").Value = 19.88
$ws.Range("V5").Value = 6.76
$ws.Range("W5").Value = 13.51
$ws.Range("X5").Value = 6.92
$ws.Range("Y5").Value = 0.94
$ws.Range("Z5").Value = 13.07
$ws.Range("AA5").Value = 5.8
$ws.Range("AB5").Value = 4.94
$ws.Range("AC5").Value = 5.86
$ws.Range("AD5").Value = 8.7
$ws.Range("AE5").Value = 0.52
$ws.Range("AF5").Value = 23.75
$ws.Range("AG5").Value = 3.6
$ws.Range("AH5").Value = 8.07

# --- Drop the last data row (row 6) entirely; dimension auto-shrinks to A1:AH5 ---
$ws.Rows("6:6").Delete()

# --- Column Z (26th column) narrows from width 8 to width 7 ---
# ColumnWidth (character units) = stored XML width - 0.83 for this workbook's font metrics.
$ws.Columns("Z").ColumnWidth = 6.17
